$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.069363594055176
$ws.Range("B1").Value = 6.332414150238037
$ws.Range("C1").Value = 5.245056629180908
$ws.Range("D1").Value = 6.097415447235107
$ws.Range("E1").Value = 4.667911529541016
